$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.072.47'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.874.77'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'312.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = "'0.5077"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.99%  '
$ws.Range("E8").Value = '  -1.80%  '
$ws.Range("D9").Value = "'0.08401"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.48%  '
$ws.Range("D11").Value = "'41.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("E12").Value = '  -2.44%  '
$ws.Range("D13").Value = '1.879.81'
$ws.Range("E13").Value = '  -0.81%  '
$ws.Range("E14").Value = '  -1.26%  '
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").Value = "'91.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").Value = "'0.06658"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").Value = "'17.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = "'6.061"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.23%  '
$ws.Range("D23").Value = '28.109.80'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = "'11.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = "'2.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("D26").Value = "'2.577"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.73%  '
$ws.Range("D27").Value = '2.097.68'
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").Value = "'157.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("D29").Value = "'20.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("D30").Value = "'125.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.73%  '
$ws.Range("D31").Value = "'0.1051"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").Value = "'1.048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").Value = "'5.624"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'3.610"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").Value = "'9.719"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.67%  '
$ws.Range("D36").Value = "'0.02453"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.33%  '
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").Value = "'0.2169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("D40").Value = "'0.6517"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("D41").Value = "'1.247"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.52%  '
$ws.Range("E42").Value = '  -0.90%  '
$ws.Range("D43").Value = "'4.897"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").Value = "'0.6174"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.21%  '
$ws.Range("D45").Value = "'13.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").Value = "'1.302"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = "'3.678"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'2.012"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("D50").Value = "'120.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("D51").Value = "'80.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.11%  '
